$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared "insert into givesOpinionConcert(" string used to build the
# SQL INSERT statements in column F is updated to the full/explicit
# "insert into givesOpinionConcert values(" form. This single shared cell
# (J1) feeds every CONCATENATE() formula in column F, so updating it here
# recalculates all the dependent cached formula values automatically.
$ws.Range("J1").Value = "insert into givesOpinionConcert values("

# Reflect where the user's selection ended up after making the edit.
$ws.Range("J2").Select() | Out-Null
